$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4770
$ws.Range("H61").Value = 229
$ws.Range("I61").Value = 50
$ws.Range("J61").Value = 318.5
$ws.Range("K61").Value = 150
$ws.Range("L61").Value = 955.5
$ws.Range("M61").Value = 22
$ws.Range("N61").Value = -1299.5
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H112").Value = 3787.1794
$ws.Range("J112").Value = 3919.4443
$ws.Range("L112").Value = 11758.3329
$ws.Range("N112").Value = -13974.3329
$ws.Range("H137").Value = 1086.4884
$ws.Range("I137").Value = 932.8919
$ws.Range("J137").Value = 2033.6666
$ws.Range("K137").Value = 2798.6757
$ws.Range("L137").Value = 6100.9998
$ws.Range("M137").Value = -248.6756999999998
$ws.Range("N137").Value = -11200.9998
$ws.Range("H138").Value = 4255.8667
$ws.Range("I138").Value = 2495.158
$ws.Range("J138").Value = 5071.8047
$ws.Range("K138").Value = 7485.474
$ws.Range("L138").Value = 15215.4141
$ws.Range("M138").Value = -2345.474
$ws.Range("N138").Value = -25495.4141
$ws.Range("H141").Value = 1084.091
$ws.Range("I141").Value = 1086.5
$ws.Range("J141").Value = 1060
$ws.Range("K141").Value = 3259.5
$ws.Range("L141").Value = 3180
$ws.Range("M141").Value = 1920.5
$ws.Range("N141").Value = -13540
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -84
$ws.Range("H32").Value = 27257.906
$ws.Range("I32").Value = 21514.049
$ws.Range("J32").Value = 145007
$ws.Range("K32").Value = 21514.049
$ws.Range("L32").Value = 145007
$ws.Range("M32").Value = -21227.049
$ws.Range("N32").Value = -145581
$ws.Range("H63").Value = 2393.2632
$ws.Range("I63").Value = 2024
$ws.Range("J63").Value = 4028.5715
$ws.Range("K63").Value = 2024
$ws.Range("L63").Value = 4028.5715
$ws.Range("M63").Value = -1338
$ws.Range("N63").Value = -5400.5715
$ws.Range("H66").Value = 2393.2632
$ws.Range("I66").Value = 2024
$ws.Range("J66").Value = 4028.5715
$ws.Range("K66").Value = 10120
$ws.Range("L66").Value = 20142.8575
$ws.Range("M66").Value = -6688
$ws.Range("N66").Value = -27006.8575
$ws.Range("H122").Value = 1721.3572
$ws.Range("J122").Value = 1399.5
$ws.Range("L122").Value = 4198.5
$ws.Range("N122").Value = -9098.5
$ws.Range("H123").Value = 29999
$ws.Range("J123").Value = 29999
$ws.Range("L123").Value = 29999
$ws.Range("N123").Value = -39799
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 39428.57
$ws.Range("J51").Value = 39428.57
$ws.Range("L51").Value = 39428.57
$ws.Range("N51").Value = -40410.57
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50546
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2617
$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1674
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H141").Value = 58052.2
$ws.Range("J141").Value = 58052.2
$ws.Range("L141").Value = 58052.2
$ws.Range("N141").Value = -68412.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3697.9443
$ws.Range("I58").Value = 820.3871
$ws.Range("J58").Value = 21538.8
$ws.Range("K58").Value = 820.3871
$ws.Range("L58").Value = 21538.8
$ws.Range("M58").Value = -617.3871
$ws.Range("N58").Value = -21944.8
$ws.Range("H122").Value = 1611.4375
$ws.Range("I122").Value = 1593.2
$ws.Range("J122").Value = 1885
$ws.Range("K122").Value = 4779.6
$ws.Range("L122").Value = 5655
$ws.Range("M122").Value = -2329.6
$ws.Range("N122").Value = -10555
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H136").Value = 3697.9443
$ws.Range("I136").Value = 820.3871
$ws.Range("J136").Value = 21538.8
$ws.Range("K136").Value = 2461.1613
$ws.Range("L136").Value = 64616.39999999999
$ws.Range("M136").Value = 88.83869999999979
$ws.Range("N136").Value = -69716.39999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H122").Value = 1476.1765
$ws.Range("I122").Value = 579.2
$ws.Range("J122").Value = 1630.8276
$ws.Range("K122").Value = 5212.8
$ws.Range("L122").Value = 14677.4484
$ws.Range("M122").Value = -2762.8
$ws.Range("N122").Value = -19577.4484
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 1098063.4
$ws.Range("I122").Value = 1197523.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3592571.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3590121.4
$ws.Range("N122").Value = -16900
$ws.Range("H123").Value = 53435.332
$ws.Range("J123").Value = 53435.332
$ws.Range("L123").Value = 53435.332
$ws.Range("N123").Value = -58335.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 2940.6155
$ws.Range("I122").Value = 2703.1875
$ws.Range("J122").Value = 3320.5
$ws.Range("K122").Value = 8109.5625
$ws.Range("L122").Value = 9961.5
$ws.Range("M122").Value = -5659.5625
$ws.Range("N122").Value = -14861.5
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 1418.1428
$ws.Range("I122").Value = 1414.3
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 4242.9
$ws.Range("L122").Value = 4485
$ws.Range("M122").Value = -1792.9
$ws.Range("N122").Value = -9385
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
